$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-25 Wednesday" "2026-02-26 Thursday"

Replace-Text "548×7=" "698×7="
Replace-Text "904×9=" "544×5="
Replace-Text "155×5=" "893×7="
Replace-Text "417×5=" "585×7="
Replace-Text "379×5=" "475×7="

Replace-Text "952×3=" "855×9="
Replace-Text "836×7=" "331×4="
Replace-Text "917×8=" "926×3="
Replace-Text "728×3=" "999×8="
Replace-Text "135×5=" "364×7="

Replace-Text "396×2=" "523×5="
Replace-Text "963×3=" "290×4="
Replace-Text "134×9=" "598×6="
Replace-Text "614×6=" "292×6="
Replace-Text "224×3=" "541×5="

Replace-Text "675×6=" "111×6="
Replace-Text "844×7=" "976×7="
Replace-Text "212×2=" "747×2="
Replace-Text "354×7=" "616×2="
Replace-Text "794×5=" "889×6="

Replace-Text "621×4=" "258×6="
Replace-Text "543×9=" "478×4="
Replace-Text "272×9=" "335×7="
Replace-Text "678×5=" "208×9="
Replace-Text "782×5=" "265×4="
